$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; existing rows 23-51 shift down to 24-52.
$ws.Rows.Item(23).Insert()

# Populate the newly-inserted row 23 with the new weekly record.
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 44915
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 100112030
$ws.Cells.Item(23, 7).Value = "Poroto granado"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 70
$ws.Cells.Item(23, 11).Value = 48000
$ws.Cells.Item(23, 12).Value = 48000
$ws.Cells.Item(23, 13).Value = 48000
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1920
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
